$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.619.54'
$ws.Range("E2").Value = '  +2.49%  '
$ws.Range("D3").Value = '2.522.92'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '2.521.91'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.35%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("E14").Value = '  +5.78%  '
$ws.Range("D15").Value = '2.982.37'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '69.541.60'
$ws.Range("E17").Value = '  +2.46%  '
$ws.Range("D18").Value = '2.499.00'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("D28").Value = '2.647.83'
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.981'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.84%  '
$ws.Range("D30").Value = '0.0₃0907'
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '508.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("E44").Value = '  -2.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.22%  '
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.518'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("D51").Value = '0.0₆0250'
$ws.Range("E51").Value = '  -2.80%  '
